# "added 4wk low sales check"
# Updates the per-week MyForecast / Inventory Coverage / Seasonality Index
# figures on the "Forecast Comparison" sheet (presumably re-computed after
# adding a rolling 4-week low-sales guard to the forecasting logic), and
# refreshes the dependent roll-up figures on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: MyForecast (D), Inventory Coverage (H), Seasonality Index (L) ---

# Week W10 (row 2)
$wsForecast.Range("D2").Value = 287
$wsForecast.Range("H2").Value = 4.12
$wsForecast.Range("L2").Value = 0.97

# Week W11 (row 3)
$wsForecast.Range("D3").Value = 293
$wsForecast.Range("H3").Value = 3.06
$wsForecast.Range("L3").Value = 1.06

# Week W12 (row 4)
$wsForecast.Range("D4").Value = 304
$wsForecast.Range("H4").Value = 1.98
$wsForecast.Range("L4").Value = 1.1

# Week W13 (row 5)
$wsForecast.Range("D5").Value = 313
$wsForecast.Range("H5").Value = 0.95
$wsForecast.Range("L5").Value = 1.06

# Week W14 (row 6)
$wsForecast.Range("D6").Value = 313
$wsForecast.Range("L6").Value = 0.87

# Week W15 (row 7)
$wsForecast.Range("D7").Value = 304
$wsForecast.Range("L7").Value = 1.1

# Week W16 (row 8)
$wsForecast.Range("D8").Value = 299
$wsForecast.Range("L8").Value = 0.82

# Week W17 (row 9)
$wsForecast.Range("D9").Value = 306
$wsForecast.Range("L9").Value = 0.87

# Week W18 (row 10)
$wsForecast.Range("D10").Value = 320
$wsForecast.Range("L10").Value = 0.91

# Week W19 (row 11)
$wsForecast.Range("D11").Value = 329
$wsForecast.Range("L11").Value = 1.19

# Week W20 (row 12)
$wsForecast.Range("D12").Value = 326
$wsForecast.Range("L12").Value = 1.13

# Week W21 (row 13)
$wsForecast.Range("D13").Value = 318
$wsForecast.Range("L13").Value = 1.18

# Week W22 (row 14)
$wsForecast.Range("D14").Value = 318
$wsForecast.Range("L14").Value = 0.96

# Week W23 (row 15)
$wsForecast.Range("D15").Value = 329
$wsForecast.Range("L15").Value = 1

# Week W24 (row 16)
$wsForecast.Range("D16").Value = 342
$wsForecast.Range("L16").Value = 0.98

# Week W25 (row 17)
$wsForecast.Range("D17").Value = 344
$wsForecast.Range("L17").Value = 0.99

# --- Summary roll-up figures (stored as text in this report) ---
$wsSummary.Range("B9").Value  = "'5052"
$wsSummary.Range("B10").Value = "'2422"
$wsSummary.Range("B11").Value = "'1199"
$wsSummary.Range("B12").Value = "'345"
$wsSummary.Range("B14").Value = "'288"
